$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: new value in column B
$ws.Range("B26").Value = 29.4

# Row 27: new value in column B and running-total formula in column C
$ws.Range("B27").Value = 23.98
$ws.Range("C27").Formula = "=B27+B26"

# Row 28: new value in column B and running-total formula in column C
$ws.Range("B28").Value = 0.13
$ws.Range("C28").Formula = "=C27+B28"

# Row 29: new value in column B and running-total formula in column C
$ws.Range("B29").Value = 1.9
$ws.Range("C29").Formula = "=C28+B29"

# Row 30: new value in column B and running-total formula in column C
$ws.Range("B30").Value = 44.57
$ws.Range("C30").Formula = "=C29+B30"

# Row 31: new value in column B and running-total formula in column C
$ws.Range("B31").Value = 0.02
$ws.Range("C31").Formula = "=C30+B31"

# Update the view: select B31 (new last cell of the data)
$ws.Range("B31").Select()
